$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '60.617.05'
$ws.Cells.Item(2, 5).Value = '  +3.25%  '

$ws.Cells.Item(3, 4).Value = '2.699.99'
$ws.Cells.Item(3, 5).Value = '  +2.93%  '

$ws.Cells.Item(4, 5).Value = '  -0.04%  '

$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '526.92'
$ws.Cells.Item(5, 5).Value = '  +1.40%  '

$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '146.27'
$ws.Cells.Item(6, 5).Value = '  +1.39%  '

$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = '0.996'
$ws.Cells.Item(7, 5).Value = '  -0.18%  '

$ws.Cells.Item(8, 5).Value = '  +1.33%  '

$ws.Cells.Item(9, 4).Value = '2.719.35'
$ws.Cells.Item(9, 5).Value = '  +3.31%  '

$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = '6.80'
$ws.Cells.Item(10, 5).Value = '  +8.27%  '

$ws.Cells.Item(11, 5).Value = '  +1.57%  '

$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = '0.340'
$ws.Cells.Item(12, 5).Value = '  +1.87%  '

$ws.Cells.Item(13, 5).Value = '  +3.25%  '

$ws.Cells.Item(14, 4).Value = '3.176.68'
$ws.Cells.Item(14, 5).Value = '  +2.87%  '

$ws.Cells.Item(15, 4).Value = '60.614.03'
$ws.Cells.Item(15, 5).Value = '  +3.23%  '

$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = '21.36'
$ws.Cells.Item(16, 5).Value = '  +2.86%  '

$ws.Cells.Item(17, 2).Value = 'ShibaInu'
$ws.Cells.Item(17, 3).Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = '0.0000138'
$ws.Cells.Item(17, 5).Value = '  +1.54%  '

$ws.Cells.Item(18, 2).Value = 'WrappedEther'
$ws.Cells.Item(18, 3).Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Cells.Item(18, 4).Value = '2.701.21'
$ws.Cells.Item(18, 5).Value = '  +2.52%  '

$ws.Cells.Item(19, 2).Value = 'Polkadot'
$ws.Cells.Item(19, 3).Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = '4.51'
$ws.Cells.Item(19, 5).Value = '  +1.23%  '

$ws.Cells.Item(20, 2).Value = 'BitcoinCash'
$ws.Cells.Item(20, 3).Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = '343.50'
$ws.Cells.Item(20, 5).Value = '  -0.36%  '

$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = '10.55'
$ws.Cells.Item(21, 5).Value = '  +3.56%  '

$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = '6.43'
$ws.Cells.Item(22, 5).Value = '  +4.74%  '

$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = '0.998'
$ws.Cells.Item(23, 5).Value = '  -0.04%  '

$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = '63.58'
$ws.Cells.Item(24, 5).Value = '  +3.62%  '

$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = '0.419'
$ws.Cells.Item(25, 5).Value = '  +1.14%  '

$ws.Cells.Item(26, 5).Value = '  +3.54%  '

$ws.Cells.Item(27, 5).Value = '  -0.23%  '

$ws.Cells.Item(28, 5).Value = '  +2.91%  '

$ws.Cells.Item(29, 5).Value = '  +3.10%  '

$ws.Cells.Item(30, 5).Value = '  +9.10%  '

$ws.Cells.Item(31, 5).Value = '  -0.12%  '

$ws.Cells.Item(32, 5).Value = '  +2.16%  '

$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = '19.00'
$ws.Cells.Item(33, 5).Value = '  +0.89%  '

$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = '149.65'
$ws.Cells.Item(34, 5).Value = '  -0.14%  '

$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = '4.27'
$ws.Cells.Item(35, 5).Value = '  +7.67%  '

$ws.Cells.Item(36, 5).Value = '  +8.21%  '

$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = '0.932'
$ws.Cells.Item(37, 5).Value = '  -3.97%  '

$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = '0.880'
$ws.Cells.Item(38, 5).Value = '  +5.39%  '

$ws.Cells.Item(39, 5).Value = '  +7.62%  '

$ws.Cells.Item(40, 5).Value = '  +1.60%  '

$ws.Cells.Item(41, 5).Value = '  +0.67%  '

$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = '281.91'
$ws.Cells.Item(42, 5).Value = '  +2.28%  '

$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = '20.18'
$ws.Cells.Item(43, 5).Value = '  +3.32%  '

$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = '0.995'
$ws.Cells.Item(44, 5).Value = '  -0.14%  '

$ws.Cells.Item(45, 5).Value = '  +1.98%  '

$ws.Cells.Item(46, 5).Value = '  +0.44%  '

$ws.Cells.Item(47, 4).Value = '2.134.93'
$ws.Cells.Item(47, 5).Value = '  +7.37%  '

$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = '4.96'
$ws.Cells.Item(48, 5).Value = '  +7.35%  '

$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = '0.0541'
$ws.Cells.Item(49, 5).Value = '  +3.76%  '

$ws.Cells.Item(50, 5).Value = '  +2.17%  '

$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = '19.26'
$ws.Cells.Item(51, 5).Value = '  +5.20%  '
